# Re-purpose the generic "2005/2007/2010/2013" workbook into a country
# performance template: sheet "2005" becomes "East Africa" (Ethiopia +
# Eritrea added), sheet "2010" becomes "Northern Europe" (Sweden added),
# and the active tab moves from the last sheet to "Northern Europe".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "2005" -> "East Africa"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "East Africa"

# Header: column B now holds "Performance" instead of the long AWU label,
# so the header row no longer needs the tall wrapped height.
$ws1.Range("B1").Value = "Performance"
$ws1.Rows.Item(1).RowHeight = 25.5

# Row 2 country: Belgium -> Ethiopia (the amount in B2 is unchanged).
$ws1.Range("A2").Value = "Ethiopia"

# New trailing row 5: Eritrea / 340. Copy row 4's formatting down first so
# the new cells pick up the same styles as the rest of the country rows.
$ws1.Range("A4").Copy($ws1.Range("A5"))
$ws1.Range("B4").Copy($ws1.Range("B5"))
$ws1.Range("A5").Value = "Eritrea"
$ws1.Range("B5").Value = 340

# Column B is now a country-name column, not the old wrapped-text column,
# so it needs to be wider (and no longer "best fit").
$ws1.Columns.Item(2).ColumnWidth = 18 - 0.833333333333

# Leave the cursor parked below the data, as in the target file.
$null = $ws1.Activate()
$ws1.Range("A6").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet 3: "2010" -> "Northern Europe"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "Northern Europe"

# Header: same "Performance" relabel + shorter row height as sheet 1.
$ws3.Range("B1").Value = "Performance"
$ws3.Rows.Item(1).RowHeight = 25.5

# Insert a new row 3 (Sweden) ahead of the existing Luxembourg/Netherlands
# rows, which shift down to rows 4 and 5 and keep their formatting.
$ws3.Rows.Item(3).Insert()
$ws3.Range("A3").Value = "Sweden"
$ws3.Range("B3").Value = 370

# Widen the country-name column here too.
$ws3.Columns.Item(2).ColumnWidth = 22 - 0.833333333333

# This sheet becomes the active tab/selection in the saved workbook.
$null = $ws3.Activate()
$ws3.Range("B11").Select() | Out-Null
